# Auto-generated script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose new values would
# otherwise be auto-converted to numbers by Excel (losing formatting,
# e.g. trailing zeros like "0.1000" -> 0.1)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.988.29"
$ws.Range("E2").Value = "  +0.25%  "

# Row 3
$ws.Range("D3").Value = "1.633.76"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").Value = "212.17"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6
$ws.Range("E6").Value = "  -0.41%  "

# Row 7
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.28%  "

# Row 8
$ws.Range("D8").Value = "23.53"
$ws.Range("E8").Value = "  -0.17%  "

# Row 9
$ws.Range("E9").Value = "  -2.16%  "

# Row 10
$ws.Range("E10").Value = "  -0.34%  "

# Row 11
$ws.Range("E11").Value = "  +0.95%  "

# Row 12
$ws.Range("D12").Value = "1.865.59"
$ws.Range("E12").Value = "  -0.56%  "

# Row 13
$ws.Range("D13").Value = "1.634.08"
$ws.Range("E13").Value = "  -0.49%  "

# Row 14
$ws.Range("E14").Value = "  -0.26%  "

# Row 15
$ws.Range("E15").Value = "  -1.77%  "

# Row 16
$ws.Range("D16").Value = "65.59"
$ws.Range("E16").Value = "  -0.29%  "

# Row 17
$ws.Range("D17").Value = "27.978.42"
$ws.Range("E17").Value = "  +0.27%  "

# Row 18
$ws.Range("D18").Value = "232.19"
$ws.Range("E18").Value = "  +0.78%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  +0.12%  "

# Row 20
$ws.Range("D20").Value = "7.56"
$ws.Range("E20").Value = "  -0.82%  "

# Row 21
$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  -0.35%  "

# Row 22
$ws.Range("E22").Value = "  -4.55%  "

# Row 23
$ws.Range("E23").Value = "  -0.78%  "

# Row 24
$ws.Range("E24").Value = "  -3.54%  "

# Row 25
$ws.Range("D25").Value = "154.57"
$ws.Range("E25").Value = "  +1.54%  "

# Row 26
$ws.Range("E26").Value = "  +0.42%  "

# Row 27
$ws.Range("E27").Value = "  -0.62%  "

# Row 28
$ws.Range("E28").Value = "  -0.37%  "

# Row 29
$ws.Range("E29").Value = "  -0.27%  "

# Row 30
$ws.Range("E30").Value = "  -0.31%  "

# Row 31
$ws.Range("E31").Value = "  -0.70%  "

# Row 32
$ws.Range("E32").Value = "  +2.06%  "

# Row 33
$ws.Range("E33").Value = "  +0.20%  "

# Row 34
$ws.Range("D34").Value = "1.408.42"
$ws.Range("E34").Value = "  -1.36%  "

# Row 35
$ws.Range("E35").Value = "  -0.16%  "

# Row 36
$ws.Range("E36").Value = "  +9.06%  "

# Row 37
$ws.Range("E37").Value = "  +0.60%  "

# Row 38
$ws.Range("E38").Value = "  +1.96%  "

# Row 39
$ws.Range("E39").Value = "  +0.17%  "

# Row 40
$ws.Range("D40").Value = "0.871"
$ws.Range("E40").Value = "  -2.07%  "

# Row 41
$ws.Range("E41").Value = "  -1.09%  "

# Row 42
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -0.25%  "

# Row 43
$ws.Range("D43").Value = "67.13"
$ws.Range("E43").Value = "  -2.16%  "

# Row 44
$ws.Range("D44").Value = "5.48"
$ws.Range("E44").Value = "  +0.66%  "

# Row 45
$ws.Range("E45").Value = "  +0.70%  "

# Row 46
$ws.Range("E46").Value = "  -0.59%  "

# Row 47
$ws.Range("D47").Value = "1.775.61"
$ws.Range("E47").Value = "  -0.49%  "

# Row 48
$ws.Range("D48").Value = "88.11"
$ws.Range("E48").Value = "  -1.08%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -3.50%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1000"
$ws.Range("E50").Value = "  -0.85%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0505"
$ws.Range("E51").Value = "  -0.25%  "

